$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update patientid (column B) values.
# Force Text format first so the purely-numeric strings are not
# auto-converted to numbers (matches original inlineStr/text cells).
$ws.Range("B2:B10").NumberFormat = "@"

$ws.Range("B2").Value = "110508"
$ws.Range("B3").Value = "33997"
$ws.Range("B4").Value = "83868"
$ws.Range("B5").Value = "92999"
$ws.Range("B6").Value = "107585"
$ws.Range("B7").Value = "133674"
$ws.Range("B8").Value = "52680"
$ws.Range("B9").Value = "58017"
$ws.Range("B10").Value = "132216"

# Fill in the previously-empty createdate (column N) values.
# Force Text format first so the date-looking strings stay as plain
# text (dd/mm/yyyy) instead of being converted to Excel date serials.
$ws.Range("N2:N10").NumberFormat = "@"

$ws.Range("N2").Value = "01/04/2025"
$ws.Range("N3").Value = "02/04/2025"
$ws.Range("N4").Value = "01/04/2025"
$ws.Range("N5").Value = "31/03/2025"
$ws.Range("N6").Value = "01/04/2025"
$ws.Range("N7").Value = "01/04/2025"
$ws.Range("N8").Value = "01/04/2025"
$ws.Range("N9").Value = "01/04/2025"
$ws.Range("N10").Value = "03/04/2025"
